$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Extend the formatting of the last existing data row (272) down to the
# new rows (273-281) so the new cells share the same "import" cell style.
$ws.Range("A272:C272").Copy()
$ws.Range("A273:C281").PasteSpecial(-4122)

$rows = @(
    @("lab.vape.menu", "Vapování"),
    @("lab.liquid.menu", "Liquidy"),
    @("lab.atomizer.menu", "Atomizéry"),
    @("lab.mod.menu", "Mody"),
    @("lab.cell.menu", "Články"),
    @("lab.vendor.menu", "Výrobci"),
    @("lab.cotton.menu", "Vaty"),
    @("lab.wire.menu", "Dráty"),
    @("lab.coil.menu", "Spirálky")
)

$startRow = 273
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $label = $rows[$i][0]
    $translation = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $label
    $ws.Cells.Item($r, 3).Value = $translation
}

$ws.Range("B277").Select() | Out-Null
